$d = $word.ActiveDocument

# --- Step 1: fix the "10/04/18" progress paragraph. It originally read
# "...de larav" + (a _GoBack bookmark) + "el et de l'option d'envoi de
# mail." - i.e. the word "laravel" was split by Word's auto _GoBack
# bookmark. A Find/Replace whose match text spans that bookmark rewrites
# the two runs into a single completed run and drops the old bookmark,
# which is exactly what the target XML shows. ---
$apos = [char]0x2019
$oldText = "larav" + "el et de l" + $apos + "option d" + $apos + "envoi de mail."
$newText = "laravel et de l" + $apos + "option d" + $apos + "envoi de mail."
$findRange = $d.Content
$found = $findRange.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

# Locate the paragraph we just rewrote, then the (empty) paragraph right
# after it - that is where the new progress note goes.
$sprintPara = $findRange.Paragraphs(1)
$targetPara = $sprintPara.Next()
$pRange = $targetPara.Range

# --- Step 2: fill in the new paragraph with two sentences, finishing with
# the _GoBack bookmark (Word always keeps exactly one, wherever the cursor
# last was). ---
$part1 = "Avancement dans la r" + [char]0x00e9 + "alisation "
$part2 = "de cr" + [char]0x00e9 + "ation de compte, s" + [char]0x00e9 + "curisation des donn" + [char]0x00e9 + "es de celui-ci."
# A trailing placeholder character is inserted too: it keeps our bookmark
# off the true paragraph end while we seat it (collapsed bookmarks placed
# exactly at a paragraph's end get mis-expanded to cover the whole
# paragraph), and is deleted again once the bookmark is in place.
$placeholder = [char]0x0001
$pRange.InsertBefore($part1 + $part2 + $placeholder)

# Re-fetch the paragraph's content range (text only, no paragraph mark).
$full = $targetPara.Range
$full.MoveEnd(1, -1)
$start = $full.Start
$len1 = $part1.Length
$len2 = $part2.Length

# Split part1/part2 into two distinct <w:r> runs: without this they would
# be re-merged into a single run (same, empty formatting) when saved. A
# zero-width bookmark dropped right on the boundary forces the split and
# is removed again afterwards, leaving the run break behind.
$midRange = $full.Duplicate
$midRange.SetRange($start + $len1, $start + $len1)
$d.Bookmarks.Add("ZZtmpSplit", $midRange) | Out-Null

# Seat the real _GoBack bookmark right after part2 (before the
# placeholder). Adding a bookmark named "_GoBack" replaces any existing
# one of that name, so the one removed in Step 1 is naturally superseded.
$endRange = $full.Duplicate
$endRange.SetRange($start + $len1 + $len2, $start + $len1 + $len2)
$d.Bookmarks.Add("_GoBack", $endRange) | Out-Null

# Drop the temporary split-point bookmark; the run split it created stays.
$d.Bookmarks("ZZtmpSplit").Delete()

# Remove the placeholder character now that _GoBack is safely seated.
$placeholderRange = $full.Duplicate
$placeholderRange.SetRange($start + $len1 + $len2, $start + $len1 + $len2 + 1)
$placeholderRange.Delete()

# --- Step 3: the document used to end with two empty paragraphs; only one
# (now filled in) is needed, so delete the paragraph mark that ends the
# paragraph we just filled in, merging it with the trailing empty one. ---
$markStart = $targetPara.Range.End - 1
$markRange = $d.Range($markStart, $markStart + 1)
$markRange.Delete() | Out-Null
